# Apply the crypto price/volume update (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'43.957.48"
$ws.Range("E2").Value = "  +3.82%  "

# Row 3
$ws.Range("D3").Value = "'2.266.65"
$ws.Range("E3").Value = "  +1.85%  "

# Row 4
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").Value = "'229.94"
$ws.Range("E5").Value = "  -0.48%  "

# Row 6
$ws.Range("D6").Value = "'0.626"
$ws.Range("E6").Value = "  +0.99%  "

# Row 7
$ws.Range("D7").Value = "'62.97"
$ws.Range("E7").Value = "  +3.47%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").Value = "'0.450"
$ws.Range("E9").Value = "  +11.54%  "

# Row 10
$ws.Range("D10").Value = "'0.101"
$ws.Range("E10").Value = "  +11.07%  "

# Row 11
$ws.Range("D11").Value = "'57.01"
$ws.Range("E11").Value = "  -0.82%  "

# Row 12
$ws.Range("B12").Value = "Avalanche"
$ws.Range("C12").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D12").Value = "'25.93"
$ws.Range("E12").Value = "  +16.33%  "

# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.106"
$ws.Range("E13").Value = "  +2.07%  "

# Row 14
$ws.Range("D14").Value = "'2.593.88"
$ws.Range("E14").Value = "  +1.37%  "

# Row 15
$ws.Range("D15").Value = "'15.64"
$ws.Range("E15").Value = "  +1.05%  "

# Row 16
$ws.Range("D16").Value = "'6.20"
$ws.Range("E16").Value = "  +9.94%  "

# Row 17
$ws.Range("D17").Value = "'0.845"
$ws.Range("E17").Value = "  +6.03%  "

# Row 18
$ws.Range("D18").Value = "'2.231.86"
$ws.Range("E18").Value = "  -0.32%  "

# Row 19
$ws.Range("D19").Value = "'43.802.98"
$ws.Range("E19").Value = "  +3.69%  "

# Row 20
$ws.Range("E20").Value = "  +5.89%  "

# Row 21
$ws.Range("D21").Value = "'73.56"
$ws.Range("E21").Value = "  +1.92%  "

# Row 22
$ws.Range("D22").Value = "'6.07"
$ws.Range("E22").Value = "  -1.64%  "

# Row 23
$ws.Range("D23").Value = "'252.96"
$ws.Range("E23").Value = "  +3.60%  "

# Row 24
$ws.Range("E24").Value = "  +0.15%  "

# Row 25
$ws.Range("E25").Value = "  -0.83%  "

# Row 26
$ws.Range("D26").Value = "'2.35"
$ws.Range("E26").Value = "  -1.31%  "

# Row 27
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'10.09"
$ws.Range("E27").Value = "  +4.19%  "

# Row 28
$ws.Range("B28").Value = "WEMIXToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D28").Value = "'3.31"
$ws.Range("E28").Value = "  +24.27%  "

# Row 29
$ws.Range("D29").Value = "'172.04"
$ws.Range("E29").Value = "  +1.75%  "

# Row 30
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'20.83"
$ws.Range("E30").Value = "  +2.42%  "

# Row 31
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Value = "'0.138"
$ws.Range("E31").Value = "  -2.28%  "

# Row 32
$ws.Range("D32").Value = "'1.39"
$ws.Range("E32").Value = "  -4.50%  "

# Row 33
$ws.Range("E33").Value = "  +3.51%  "

# Row 34
$ws.Range("D34").Value = "'0.0695"
$ws.Range("E34").Value = "  +6.63%  "

# Row 35
$ws.Range("D35").Value = "'4.78"
$ws.Range("E35").Value = "  +0.88%  "

# Row 36
$ws.Range("D36").Value = "'4.91"
$ws.Range("E36").Value = "  -1.30%  "

# Row 37
$ws.Range("E37").Value = "  +7.37%  "

# Row 38
$ws.Range("D38").Value = "'6.52"
$ws.Range("E38").Value = "  +2.33%  "

# Row 39
$ws.Range("E39").Value = "  -1.77%  "

# Row 40
$ws.Range("D40").Value = "'0.0258"
$ws.Range("E40").Value = "  +3.57%  "

# Row 41
$ws.Range("E41").Value = "  -0.27%  "

# Row 42
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'17.57"
$ws.Range("E42").Value = "  +9.07%  "

# Row 43
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "'0.0974"
$ws.Range("E43").Value = "  +1.10%  "

# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'8.22"
$ws.Range("E44").Value = "  -5.03%  "

# Row 45
$ws.Range("B45").Value = "TerraClassic"
$ws.Range("C45").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D45").Value = "'0.000216"
$ws.Range("E45").Value = "  -3.98%  "

# Row 46
$ws.Range("D46").Value = "'98.49"
$ws.Range("E46").Value = "  +1.64%  "

# Row 47
$ws.Range("E47").Value = "  -0.40%  "

# Row 48
$ws.Range("B48").Value = "Celestia"
$ws.Range("C48").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D48").Value = "'10.06"
$ws.Range("E48").Value = "  +16.37%  "

# Row 49
$ws.Range("D49").Value = "'1.450.41"
$ws.Range("E49").Value = "  -0.30%  "

# Row 50
$ws.Range("B50").Value = "FTXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D50").Value = "'4.30"
$ws.Range("E50").Value = "  -1.08%  "

# Row 51
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'2.29"
$ws.Range("E51").Value = "  +3.97%  "
